# Update "Diferença" report from 2023/03-2022/03 to 2023/04-2022/04.
# New month's ranking reorders/changes the regions and adds one more row
# (Amazonas / Rio Grande do Norte / Rio de Janeiro enter the table,
#  Distrito Federal and Bahia drop out), so the data block grows from
# 8 rows (A2:D9) to 9 rows (A2:D10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$label = "Diferença 2023/04 - 2022/04"

# region, value, placement (placement empty string means no "D" entry)
$data = @(
    @("Acre",                 3.602327346319242,  "1º"),
    @("Espírito Santo",       1.963631642247975,  "2º"),
    @("Rio Grande do Norte",  1.606708664898335,  "3º"),
    @("Rio de Janeiro",       1.397378335914638,  "4º"),
    @("Maranhão",             1.236047897535343,  "5º"),
    @("Amazonas",             1.206775614614415,  "6º"),
    @("Sergipe",              0.6547022702024776, "11º"),
    @("Nordeste",             0.4249213300867609, ""),
    @("Brasil",               0.5320266924128561, "")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $label
    $ws.Cells.Item($row, 3).Value = $item[1]
    if ($item[2] -ne "") {
        $ws.Cells.Item($row, 4).Value = $item[2]
    } else {
        $ws.Cells.Item($row, 4).Value = $null
    }
    $row++
}
